# edit.ps1 -- apply the "fechaEliminacion" / "apellido" placeholder
# refactor described by the commit diff.
#
# Word COM constants used below:
#   wdFindContinue = 1
#   wdReplaceAll   = 2
#   wdReplaceNone  = 0

$d = $word.ActiveDocument

function Merge-Range([int]$start, [int]$end) {
    # Forces Word to re-flow / coalesce the run(s) spanning [start,end)
    # into a single run without altering the visible text: append a
    # one-character sentinel right after the span (inheriting the
    # formatting of the span's own tail), then delete that sentinel.
    $full = $d.Range($start, $end)
    $full.InsertAfter("Z")
    $sentinel = $d.Range($end, $end + 1)
    $sentinel.Text = ""
}

# ---------------------------------------------------------------------
# 1) "Medellín, 27 de abril 2024. Señor (a)"  ->  "Medellín, {fechaEliminacion}. Señor (a)"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("27 de abril 2024", $false, $false, $false, $false, $false, `
    $true, 1, $false, "{fechaEliminacion}", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Heading paragraph "{nombre}"  ->  "{nombre} {apellido}"
# ---------------------------------------------------------------------
$pNombreHeading = $d.Paragraphs.Item(3).Range
$pNombreHeading.Find.Execute("{nombre}", $false, $false, $false, $false, $false, `
    $true, 1, $false, "{nombre} {apellido}", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) ASUNTO paragraph: "... (a) {nombre} el {fechaFinal}"
#    -> "... (a) {nombre} {apellido} el {fechaEliminacion}"
# ---------------------------------------------------------------------
$pAsunto = $d.Paragraphs.Item(8).Range
$pAsunto.Find.Execute("{nombre}", $false, $false, $false, $false, $false, `
    $true, 1, $false, "{nombre} {apellido}", 2) | Out-Null

$pAsunto2 = $d.Paragraphs.Item(8).Range
$pAsunto2.Find.Execute("{fechaFinal}", $false, $false, $false, $false, $false, `
    $true, 1, $false, "{fechaEliminacion}", 2) | Out-Null

# ---------------------------------------------------------------------
# 4) "Esta aceptación es efectiva a partir del día {fechaFinal}, fecha..."
#    -> "...del día {fechaEliminacion}, fecha..." (also coalesces the
#    previously-split "dí" / "a " runs into a single "día " run)
# ---------------------------------------------------------------------
$pDia = $d.Paragraphs.Item(13).Range
$pDia.Find.Execute("{fechaFinal}", $false, $false, $false, $false, $false, `
    $true, 1, $false, "{fechaEliminacion}", 2) | Out-Null

$pDia2 = $d.Paragraphs.Item(13).Range
$rDia = $pDia2.Duplicate
$rDia.Find.Execute("dí", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
Merge-Range $rDia.Start ($rDia.End + 2)

# ---------------------------------------------------------------------
# 5) "...en atención a q" + "ue esta ha sido..." -> single run (no text
#    change, just coalesces the two plain runs without touching the
#    surrounding bold "ELEVEN TWO S.A.S" runs)
# ---------------------------------------------------------------------
$pElev = $d.Paragraphs.Item(15).Range
$r1 = $pElev.Duplicate
$r1.Find.Execute("acepta su renuncia tomando en cuenta", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$r2 = $pElev.Duplicate
$r2.Find.Execute("y en ningún momento empresa ", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
Merge-Range $r1.Start $r2.End

# ---------------------------------------------------------------------
# 6) "...presentar s" + "u renuncia, le agradecemos" -> single run
# ---------------------------------------------------------------------
$pAsi = $d.Paragraphs.Item(17).Range
$r3 = $pAsi.Duplicate
$r3.Find.Execute("Así las cosas, respetamos los motivos personales que tenga para presentar s", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r4 = $pAsi.Duplicate
$r4.Find.Execute("u renuncia, le agradecemos", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
Merge-Range $r3.Start $r4.End

# ---------------------------------------------------------------------
# 7) "...mejoramiento continuo de " + "nuestros procesos." -> single run
# ---------------------------------------------------------------------
$pProc = $d.Paragraphs.Item(17).Range
$r5 = $pProc.Duplicate
$r5.Find.Execute("conocimientos que compartió con esta empresa, se traducirán en un mejoramiento continuo de ", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r6 = $pProc.Duplicate
$r6.Find.Execute("nuestros procesos.", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
Merge-Range $r5.Start $r6.End

# ---------------------------------------------------------------------
# 8) Signature table cell: "{nombre }"  ->  "{nombre} {apellido}"
# ---------------------------------------------------------------------
$pTableName = $d.Paragraphs.Item(38).Range
$pTableName.Find.Execute("{nombre }", $false, $false, $false, $false, $false, `
    $true, 1, $false, "{nombre} {apellido}", 2) | Out-Null
